$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update input values; dependent formula cells recalc automatically
$ws.Range("A5").Value = 1572.24
$ws.Range("C8").Value = 2347.5

# Recalculate the workbook to refresh formula results
$excel.Calculate()

# Update the selected cell / active cell to C5
$ws.Range("C5").Select()
